# Rename the sheet "Property1" to "DataNode" (unify DataNode/DataTable/Entity naming)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Property1")
$ws.Name = "DataNode"

# Move the active selection to H33, matching the saved selection state in the diff
$ws.Range("H33").Select()
